$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Capture existing data before overwriting anything -----------------
# Old layout: A = segment label (bold/border/center style), B = PercActivations,
#             C = PercSegmentAreas. Row 1 only has B1/C1 headers.
$labels = @()
$perc = @()
$area = @()
for ($r = 2; $r -le 20; $r++) {
    $labels += , $ws.Cells.Item($r, 1).Value2
    $perc += , $ws.Cells.Item($r, 2).Value2
    $area += , $ws.Cells.Item($r, 3).Value2
}
$header1 = $ws.Cells.Item(1, 2).Value2   # "PercActivations"
$header2 = $ws.Cells.Item(1, 3).Value2   # "PercSegmentAreas"

# Grab the header/label style (bold, centered, bordered - cellXf "1") before
# the layout shuffle so every relocated header can reuse the same style.
$ws.Range("B1").Copy()

# --- Shift the two data columns one place to the right ------------------
# New layout: A = numeric segment index, B = segment label, C = PercActivations,
#             D = PercSegmentAreas.
$ws.Cells.Item(1, 4).Value = $header2
$ws.Cells.Item(1, 3).Value = $header1
$ws.Cells.Item(1, 2).Value = "segments"

$ws.Range("B1:D1").PasteSpecial($xlPasteFormats)

for ($r = 2; $r -le 20; $r++) {
    $i = $r - 2

    # Column A keeps the bold/border/centered formatting, but now holds the
    # numeric 0-based segment index instead of the label text.
    $ws.Cells.Item($r, 1).Value = $i

    # Column B gets the plain (unstyled) label text that used to live in A.
    $ws.Cells.Item($r, 2).Value = $labels[$i]
    $ws.Cells.Item($r, 2).ClearFormats()

    # Columns C and D get the data values, shifted over from B and C.
    $ws.Cells.Item($r, 3).Value = $perc[$i]
    $ws.Cells.Item($r, 4).Value = $area[$i]
}

$wb.Save()
